# PV_PROPOSAL_PROJECT_TyperZombies.docx edits
# Commit message: "posisi y spawn + game over"
#
# Strategy: locate the exact target paragraph for each change (by scanning
# Paragraphs for distinctive text), compute a precise sub-Range covering only
# the text that needs restructuring, and use Range.InsertXML with an explicit
# run-level payload so the resulting <w:r> boundaries match the target
# exactly (InsertXML replaces only the content of the Range it's called on).

$d = $word.ActiveDocument

function New-RunXml($text, $preserve) {
    if ($preserve) {
        $space = ' xml:space="preserve"'
    } else {
        $space = ""
    }
    return '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
}

function Wrap-Package($innerRuns) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerRuns + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphContaining($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

function Replace-SubRange($para, $oldText, $innerRuns) {
    $full = $para.Range.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "text not found: $oldText"
    }
    $s = $para.Range.Start + $idx
    $e = $s + $oldText.Length
    $sub = $d.Range($s, $e)
    $xml = Wrap-Package $innerRuns
    $sub.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from its old spot (Highscore paragraph,
#    after "gold"). It is re-added below near "dibawah".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "... score - gold - level ..." : merge " - " + "gold" runs into one
#    run " - gold" (bookmark already removed above).
# ---------------------------------------------------------------------
$p1 = Find-ParagraphContaining "score - gold - level"
$runs1 = New-RunXml " - gold" $true
Replace-SubRange $p1 " - gold" $runs1

# ---------------------------------------------------------------------
# 3) "... atas. Isi menu adalah:" -> "... atas. Isi menu" + " pause" + " adalah:"
# ---------------------------------------------------------------------
$p2 = Find-ParagraphContaining "Isi menu adalah:"
$runs2 = (New-RunXml ". Isi menu" $false) + (New-RunXml " pause" $true) + (New-RunXml " adalah:" $true)
Replace-SubRange $p2 ". Isi menu adalah:" $runs2

# ---------------------------------------------------------------------
# 4) "... sesuai dengan tabel dibawah." -> "... dibawah" + "." + bookmark
# ---------------------------------------------------------------------
$p3 = Find-ParagraphContaining "sesuai dengan tabel dibawah"
$runs3 = (New-RunXml " sesuai dengan tabel dibawah" $true) + (New-RunXml "." $false) + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Replace-SubRange $p3 " sesuai dengan tabel dibawah." $runs3

# ---------------------------------------------------------------------
# 5) "...sampai 3 setiap 3 detik." -> "...sampai 2 setiap 3 detik."
#    ("3 setiap " run splits into "2" + " setiap " runs)
# ---------------------------------------------------------------------
$p4 = Find-ParagraphContaining "sampai 3 setiap"
$runs4 = (New-RunXml "2" $false) + (New-RunXml " setiap " $true)
Replace-SubRange $p4 "3 setiap " $runs4

# ---------------------------------------------------------------------
# 6) Game over paragraph: split the trailing sentence and add new text
#    about returning to main menu / new game / continue.
# ---------------------------------------------------------------------
$p5 = Find-ParagraphContaining "kembali ke level 1 dan status nya akan di reset pula"
$oldTail = ". Bila player kehabisan hp, maka akan ditampilkan pesan game over dan player akan kembali ke level 1 dan status nya akan di reset pula."
$runs5 = (New-RunXml ". Bila player kehabisan hp, maka akan ditampilkan pesan game over " $true) `
    + (New-RunXml "dan player akan kembali ke main menu" $false) `
    + (New-RunXml "." $false) `
    + (New-RunXml " Player yang sudah game over dapat new game menggunakan nama yang sama" $true) `
    + (New-RunXml ", tetapi tidak dapat melakukan continue" $false) `
    + (New-RunXml "." $false)
Replace-SubRange $p5 $oldTail $runs5

Write-Output "Done"
